# Mark "Usage Status" (column L) as "Y" for calibration rows that were
# previously blank (or, for rows 26-28, marked "N"). This mirrors the
# author's commit "Update Excel file w/ usage status of calibration data".
#
# Rows 123:138 and 155:170 additionally need their fill color switched to
# match the lighter "Usage Status" shading used elsewhere in the table
# (style index 18's fill, FFEBFAFF) instead of the darker shading they
# inherited from the surrounding J/K/M cells (style index 20's fill,
# FFC1EFFF).
#
# NB: multi-area ("L1:L2,L4:L5") Range.Value assignment only applies to
# the first area in this engine, so each contiguous block is set with its
# own statement instead of being combined into one Union/comma range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows whose L cell is simply filled in with "Y" (no style change) ---
$ws.Range("L10:L22").Value = "Y"
$ws.Range("L26:L40").Value = "Y"
$ws.Range("L97:L122").Value = "Y"
$ws.Range("L139:L154").Value = "Y"
$ws.Range("L171:L174").Value = "Y"

# --- Rows whose L cell gets "Y" plus a fill-color fix ---
$ws.Range("L123:L138").Value = "Y"
$ws.Range("L123:L138").Interior.Color = 16775915

$ws.Range("L155:L170").Value = "Y"
$ws.Range("L155:L170").Interior.Color = 16775915

# Leave the view scrolled/selected roughly where the author left it.
[void]$ws.Range("F97").Select()
